# User data 3.0 to dev
# The "Data-wide-value" sheet gains a new "budget-type" column (inserted as
# column B), pushing the existing year columns (2012-2016) one slot to the
# right (B:F -> C:G). Every data row is stamped with the constant value
# "budget" in the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Insert a new blank column before the first year column (old column B).
$ws.Range("B1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("B1").Value = "budget-type"

# Every data row (2 through 112) gets the constant "budget" value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 112
}
$ws.Range("B2:B" + $lastRow).Value = "budget"
